$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new value in A3
$ws.Range("A3").Value = 3183738

# Update the selected cell to match the final state (B7)
$ws.Range("B7").Select()
